$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.457.66"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").Value = "1.571.44"
$ws.Range("E3").Value = "  +0.45%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D6").Value = "'288.61"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("D7").Value = "'0.3707"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.91%  "

$ws.Range("D8").Value = "'48.34"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.06%  "

$ws.Range("D9").Value = "'0.3309"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.64%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.07503"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.57%  "

$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "'1.135"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.00%  "

$ws.Range("E12").Value = "  -0.04%  "

$ws.Range("D13").Value = "'20.74"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.89%  "

$ws.Range("D14").Value = "'5.932"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.76%  "

$ws.Range("D15").Value = "'6.871"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.86%  "

$ws.Range("D16").Value = "1.569.15"
$ws.Range("E16").Value = "  +0.32%  "

$ws.Range("E17").Value = "  +0.90%  "

$ws.Range("E18").Value = "  +0.24%  "

$ws.Range("D19").Value = "'87.62"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.57%  "

$ws.Range("E20").Value = "  -0.02%  "

$ws.Range("D21").Value = "'6.349"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.33%  "

$ws.Range("D22").Value = "'16.53"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.25%  "

$ws.Range("D23").Value = "'12.04"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").Value = "22.445.18"
$ws.Range("E24").Value = "  +0.11%  "

$ws.Range("D25").Value = "'2.392"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("D26").Value = "'2.584"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.03%  "

$ws.Range("D27").Value = "'153.66"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.00%  "

$ws.Range("D28").Value = "'19.69"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.37%  "

$ws.Range("D29").Value = "'5.010"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.63%  "

$ws.Range("D30").Value = "'124.47"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.48%  "

$ws.Range("D31").Value = "1.745.31"
$ws.Range("E31").Value = "  +0.45%  "

$ws.Range("D32").Value = "'1.059"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.73%  "

$ws.Range("E33").Value = "  +0.24%  "

$ws.Range("D34").Value = "'6.117"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.21%  "

$ws.Range("D35").Value = "'9.796"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.11%  "

$ws.Range("D36").Value = "'0.08374"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.20%  "

$ws.Range("D37").Value = "'0.02467"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.25%  "

$ws.Range("D38").Value = "'0.2263"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.08%  "

$ws.Range("D39").Value = "'0.06409"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.13%  "

$ws.Range("D40").Value = "'1.287"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.42%  "

$ws.Range("D41").Value = "'5.339"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.26%  "

$ws.Range("D42").Value = "'0.6317"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.96%  "

$ws.Range("E43").Value = "  +1.69%  "

$ws.Range("D44").Value = "'13.91"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.84%  "

$ws.Range("D45").Value = "'0.6175"
$ws.Range("D45").ClearFormats()

$ws.Range("E46").Value = "  +0.26%  "

$ws.Range("D47").Value = "'2.062"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.75%  "

$ws.Range("D48").Value = "'125.87"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.27%  "

$ws.Range("E49").Value = "  -0.55%  "

$ws.Range("D50").Value = "'0.07221"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.48%  "

$ws.Range("D51").Value = "'76.85"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.38%  "

